# Update cryptocurrency price/volume data per Sep 2, 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.006.26"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.450.80"
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.09"
$ws.Range("E5").Value = "  -1.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.17"
$ws.Range("E6").Value = "  -1.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").Value = "  -1.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.473.54"
$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("E10").Value = "  -3.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("E12").Value = "  -2.54%  "

$ws.Range("E13").Value = "  -3.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.889.79"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.962.16"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.02"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.465.45"
$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.56"
$ws.Range("E19").Value = "  -2.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "319.08"
$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("E21").Value = "  -1.32%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.98"
$ws.Range("E22").Value = "  +3.73%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.47"
$ws.Range("E24").Value = "  -1.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.405"
$ws.Range("E25").Value = "  -1.45%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.23%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("E28").Value = "  -1.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.96"
$ws.Range("E29").Value = "  +0.85%  "

$ws.Range("E30").Value = "  -3.12%  "

$ws.Range("E31").Value = "  -2.36%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.19"
$ws.Range("E32").Value = "  -2.44%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("E33").Value = "  -1.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.83"
$ws.Range("E36").Value = "  -1.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("E37").Value = "  -3.57%  "

$ws.Range("E38").Value = "  -1.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.59"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  -1.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.761"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "272.76"
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("E43").Value = "  -2.78%  "

$ws.Range("E44").Value = "  -0.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.588"
$ws.Range("E45").Value = "  -1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.92"
$ws.Range("E47").Value = "  -5.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0491"
$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.35"
$ws.Range("E49").Value = "  -3.16%  "

$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.72"
$ws.Range("E51").Value = "  -2.70%  "
